$wb = $excel.ActiveWorkbook

# The localized report moved from "Ready for handoff" to "In Translation" for
# every locale status cell (Overview!E2:F4 plus the per-locale Status column
# on the "zh-cn" and "de-de" sheets). Shrinking that text means the
# auto-fitted Status/locale columns get narrower too, so re-fit them after
# the value change.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"
$zhcn.Columns.Item(3).AutoFit() | Out-Null
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"
$dede.Columns.Item(3).AutoFit() | Out-Null
$dede.Columns.Item(3).ColumnWidth = 12.5
